$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 21.01.2022 15:30"

# Row 2: new price check result.
# The previous "current price" (B2) becomes the "old price" (C2),
# and the freshly scraped price becomes the new B2.
$ws.Range("C2").Value = $ws.Range("B2").Value2
$ws.Range("B2").Value = 34.5

# Delta is now written as a text string instead of a number.
# Mark the cell as Text first so Excel doesn't auto-convert "+0.6" to a
# number, then reset the style back to the default (no custom formatting).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "+0.6"
$ws.Range("D2").Style = "Normal"

# Date/time is now written as plain text (not a numeric date serial).
# Same trick: force Text format so the string isn't parsed as a date,
# then reset the style back to the default.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2022-01-21 15:30:04"
$ws.Range("E2").Style = "Normal"
